$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "correlation of determination" (R^2) spike-solution cells.
$ws.Range("I8").Formula = "=SUM(sum_x*sum_y)"
$ws.Range("C8").Copy()
$ws.Range("I8").PasteSpecial(-4122)  # xlPasteFormats - match the yellow "result" style used by row 8

$ws.Range("L9").Formula = "=SUM(5*74498)"
$ws.Range("M20").Formula = "=SUM(5*74498)"
$ws.Range("M21").Formula = "=SUM(580*580)"
$ws.Range("M22").Formula = "=M20-M21"

# View state: zoom back to 100% and move the selection/scroll position.
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("M23").Select()
